# Command class: make non-abstract
#
# "Rectangle 45" on slide 2 is the {abstract}/Command class box in the
# diagram. Its text is currently two runs separated by a soft line break:
#   "{abstract}" <br/> "Command"
# We want it to read just "Command" (dropping the "{abstract}" marker and
# the now-unneeded line break), so students who haven't covered abstract
# classes aren't confused by it.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item("Rectangle 45")
$tr = $sh.TextFrame.TextRange

# Original run layout (1-based character positions):
#   1-10  "{abstract}"
#   11    line break (<a:br/>)
#   12-18 "Command"
#
# First remove characters 10-18 (the last letter of "{abstract}" together
# with the line break and the trailing "Command" run) -- spanning across
# the <a:br/> like this is what actually drops the line-break node, not
# just the text around it.
$tail = $tr.Characters(10, 9)
$tail.Text = ""

# That left us with "{abstrac" (the trailing "t" was removed above); put
# back the final word by replacing the remaining "{abstrac" text with
# "Command".
$head = $tr.Characters(1, 9)
$head.Text = "Command"

# Best-effort: the source deck also gained two centered slide guides
# (horizontal @ 270pt, vertical @ 360pt) as an incidental editor-state
# change. Add them through the standard Guides collection; harmless no-op
# if the host doesn't persist presentation-level guides.
try {
    [void]$p.Guides.Add(1, 270)
    [void]$p.Guides.Add(2, 360)
} catch {
}
